$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value2 = 3500
$ws.Range("J13").Value2 = 3500
$ws.Range("L13").Value2 = 3500
$ws.Range("N13").Value2 = -3838
# Row 17
$ws.Range("H17").Value2 = 1428.3726
$ws.Range("I17").Value2 = 705
$ws.Range("J17").Value2 = 1604.8049
$ws.Range("K17").Value2 = 2115
$ws.Range("L17").Value2 = 4814.4147
$ws.Range("M17").Value2 = -1947
$ws.Range("N17").Value2 = -5150.4147
# Row 40
$ws.Range("H40").Value2 = 4172041.2
$ws.Range("I40").Value2 = 3166.3333
$ws.Range("K40").Value2 = 3166.3333
$ws.Range("M40").Value2 = -2991.3333
# Row 137
$ws.Range("H137").Value2 = 4699.7144
$ws.Range("I137").Value2 = 8940.200000000001
$ws.Range("K137").Value2 = 26820.6
$ws.Range("M137").Value2 = -24270.6
# Row 138
$ws.Range("H138").Value2 = 1670951.9
$ws.Range("I138").Value2 = 1434.2307
$ws.Range("K138").Value2 = 4302.6921
$ws.Range("M138").Value2 = 837.3078999999998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 1928769.6
$ws.Range("I32").Value2 = 1988809.6
$ws.Range("K32").Value2 = 1988809.6
$ws.Range("M32").Value2 = -1988522.6
# Row 45
$ws.Range("H45").Value2 = 4346.5557
$ws.Range("I45").Value2 = 1818.2
$ws.Range("K45").Value2 = 1818.2
$ws.Range("M45").Value2 = -1441.2
# Row 74
$ws.Range("H74").Value2 = 18473.742
$ws.Range("I74").Value2 = 23058.957
$ws.Range("K74").Value2 = 23058.957
$ws.Range("M74").Value2 = -22184.957
# Row 77
$ws.Range("H77").Value2 = 18473.742
$ws.Range("I77").Value2 = 23058.957
$ws.Range("K77").Value2 = 115294.785
$ws.Range("M77").Value2 = -110926.785
# Row 132
$ws.Range("H132").Value2 = 3618.87
$ws.Range("I132").Value2 = 1735.9122
$ws.Range("K132").Value2 = 5207.7366
$ws.Range("M132").Value2 = -2677.7366

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value2 = 1554.2916
$ws.Range("I94").Value2 = 1073
$ws.Range("J94").Value2 = 2356.4443
$ws.Range("K94").Value2 = 1073
$ws.Range("L94").Value2 = 2356.4443
$ws.Range("M94").Value2 = -622
$ws.Range("N94").Value2 = -3258.4443
# Row 95
$ws.Range("H95").Value2 = 59250
$ws.Range("J95").Value2 = 59250
$ws.Range("L95").Value2 = 59250
$ws.Range("N95").Value2 = -64742
# Row 134
$ws.Range("H134").Value2 = 4063.0137
$ws.Range("I134").Value2 = 1470.3019
$ws.Range("K134").Value2 = 4410.905699999999
$ws.Range("M134").Value2 = -1875.905699999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value2 = 1865
$ws.Range("J10").Value2 = 2137.6
$ws.Range("L10").Value2 = 2137.6
$ws.Range("N10").Value2 = -2415.6
# Row 31
$ws.Range("H31").Value2 = 7264.625
$ws.Range("I31").Value2 = 3804.25
$ws.Range("K31").Value2 = 3804.25
$ws.Range("M31").Value2 = -3509.25
# Row 34
$ws.Range("H34").Value2 = 7264.625
$ws.Range("I34").Value2 = 3804.25
$ws.Range("K34").Value2 = 3804.25
$ws.Range("M34").Value2 = -3602.25
# Row 94
$ws.Range("H94").Value2 = 1165.8636
$ws.Range("I94").Value2 = 1608.3334
$ws.Range("J94").Value2 = 999.9375
$ws.Range("K94").Value2 = 1608.3334
$ws.Range("L94").Value2 = 999.9375
$ws.Range("M94").Value2 = -1157.3334
$ws.Range("N94").Value2 = -1901.9375
# Row 119
$ws.Range("H119").Value2 = 0
$ws.Range("J119").Value2 = 0
$ws.Range("L119").Value2 = 0
$ws.Range("N119").ClearContents()
# Row 134
$ws.Range("H134").Value2 = 2316.97
$ws.Range("I134").Value2 = 1154.5143
$ws.Range("J134").Value2 = 5029.3667
$ws.Range("K134").Value2 = 3463.5429
$ws.Range("L134").Value2 = 15088.1001
$ws.Range("M134").Value2 = -928.5429000000004
$ws.Range("N134").Value2 = -20158.1001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 2862.818
$ws.Range("I5").Value2 = 988.75
$ws.Range("K5").Value2 = 2966.25
$ws.Range("M5").Value2 = -2854.25
# Row 16
$ws.Range("H16").Value2 = 3400
$ws.Range("J16").Value2 = 3400
$ws.Range("L16").Value2 = 10200
$ws.Range("N16").Value2 = -10546
# Row 75
$ws.Range("H75").Value2 = 31749334
$ws.Range("I75").Value2 = 83333544
$ws.Range("K75").Value2 = 250000632
$ws.Range("M75").Value2 = -249999634
# Row 78
$ws.Range("H78").Value2 = 31749334
$ws.Range("I78").Value2 = 83333544
$ws.Range("K78").Value2 = 750001896
$ws.Range("M78").Value2 = -749996904
# Row 88
$ws.Range("H88").Value2 = 3000
$ws.Range("J88").Value2 = 3000
$ws.Range("L88").Value2 = 9000
$ws.Range("N88").Value2 = -9856
# Row 91
$ws.Range("H91").Value2 = 3000
$ws.Range("J91").Value2 = 3000
$ws.Range("L91").Value2 = 9000
$ws.Range("N91").Value2 = -11964
# Row 135
$ws.Range("H135").Value2 = 2862.818
$ws.Range("I135").Value2 = 988.75
$ws.Range("K135").Value2 = 8898.75
$ws.Range("M135").Value2 = -6363.75
# Row 137
$ws.Range("H137").Value2 = 127680.94
$ws.Range("J137").Value2 = 146771.14
$ws.Range("L137").Value2 = 440313.42
$ws.Range("N137").Value2 = -450513.42

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value2 = 53942
$ws.Range("J39").Value2 = 53942
$ws.Range("L39").Value2 = 53942
$ws.Range("N39").Value2 = -55006
# Row 80
$ws.Range("H80").Value2 = 3254.5386
$ws.Range("I80").Value2 = 2946.1428
$ws.Range("J80").Value2 = 3614.3333
$ws.Range("K80").Value2 = 2946.1428
$ws.Range("L80").Value2 = 3614.3333
$ws.Range("M80").Value2 = -1948.1428
$ws.Range("N80").Value2 = -5610.3333
# Row 83
$ws.Range("H83").Value2 = 3254.5386
$ws.Range("I83").Value2 = 2946.1428
$ws.Range("J83").Value2 = 3614.3333
$ws.Range("K83").Value2 = 14730.714
$ws.Range("L83").Value2 = 18071.6665
$ws.Range("M83").Value2 = -9738.714
$ws.Range("N83").Value2 = -28055.6665
# Row 113
$ws.Range("H113").Value2 = 5624.093
$ws.Range("I113").Value2 = 2930.5
$ws.Range("K113").Value2 = 2930.5
$ws.Range("M113").Value2 = -760.5
# Row 132
$ws.Range("H132").Value2 = 4025.087
$ws.Range("I132").Value2 = 1407.4
$ws.Range("J132").Value2 = 8933.25
$ws.Range("K132").Value2 = 4222.200000000001
$ws.Range("L132").Value2 = 26799.75
$ws.Range("M132").Value2 = -1692.200000000001
$ws.Range("N132").Value2 = -31859.75
# Row 140
$ws.Range("H140").Value2 = 70690
$ws.Range("J140").Value2 = 70690
$ws.Range("L140").Value2 = 70690
$ws.Range("N140").Value2 = -81050
# Row 141
$ws.Range("H141").Value2 = 66244
$ws.Range("I141").Value2 = 64988
$ws.Range("J141").Value2 = 67500
$ws.Range("K141").Value2 = 64988
$ws.Range("L141").Value2 = 67500
$ws.Range("M141").Value2 = -59808
$ws.Range("N141").Value2 = -77860

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value2 = 5818155.5
$ws.Range("I132").Value2 = 8335322.5
$ws.Range("J132").Value2 = 9307.576999999999
$ws.Range("K132").Value2 = 25005967.5
$ws.Range("L132").Value2 = 27922.731
$ws.Range("M132").Value2 = -25003437.5
$ws.Range("N132").Value2 = -32982.731
# Row 136
$ws.Range("H136").Value2 = 7619.4893
$ws.Range("I136").Value2 = 2504.625
$ws.Range("J136").Value2 = 12956.739
$ws.Range("K136").Value2 = 7513.875
$ws.Range("L136").Value2 = 38870.217
$ws.Range("M136").Value2 = -4963.875
$ws.Range("N136").Value2 = -43970.217

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value2 = 107971.13
$ws.Range("I122").Value2 = 135233.77
$ws.Range("J122").Value2 = 5736.25
$ws.Range("K122").Value2 = 405701.3099999999
$ws.Range("L122").Value2 = 17208.75
$ws.Range("M122").Value2 = -403251.3099999999
$ws.Range("N122").Value2 = -22108.75
# Row 132
$ws.Range("H132").Value2 = 16952466
$ws.Range("I132").Value2 = 19233818
$ws.Range("J132").Value2 = 5286.143
$ws.Range("K132").Value2 = 57701454
$ws.Range("L132").Value2 = 15858.429
$ws.Range("M132").Value2 = -57698924
$ws.Range("N132").Value2 = -20918.429
# Row 141
$ws.Range("H141").Value2 = 100000
$ws.Range("J141").Value2 = 100000
$ws.Range("L141").Value2 = 100000
$ws.Range("N141").Value2 = -110360
